$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.18
$ws.Range("C2").Value = 0.5733333333333334
$ws.Range("J2").Value = 0.01333333333333333
$ws.Range("P2").Value = 0.13
$ws.Range("S2").Value = 0.1033333333333333
$ws.Range("B3").Value = 0.01675977653631285
$ws.Range("C3").Value = 0.0223463687150838
$ws.Range("J3").Value = 0.02793296089385475
$ws.Range("P3").Value = 0.7430167597765364
$ws.Range("S3").Value = 0.1899441340782123
$ws.Range("J4").Value = 0.07692307692307693
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.2307692307692308
$ws.Range("B6").Value = 0.04329004329004329
$ws.Range("D6").Value = 0.01298701298701299
$ws.Range("F6").Value = 0.08658008658008658
$ws.Range("J6").Value = 0.2857142857142857
$ws.Range("O6").Value = 0.008658008658008658
$ws.Range("Q6").Value = 0.2164502164502164
$ws.Range("R6").Value = 0.04761904761904762
$ws.Range("S6").Value = 0.2987012987012987
$ws.Range("B7").Value = 0.08203125
$ws.Range("D7").Value = 0.01953125
$ws.Range("F7").Value = 0.04296875
$ws.Range("J7").Value = 0.16796875
$ws.Range("O7").Value = 0.01953125
$ws.Range("Q7").Value = 0.203125
$ws.Range("R7").Value = 0.046875
$ws.Range("S7").Value = 0.41796875
$ws.Range("B8").Value = 0.1043256997455471
$ws.Range("D8").Value = 0.01272264631043257
$ws.Range("F8").Value = 0.06361323155216285
$ws.Range("J8").Value = 0.1450381679389313
$ws.Range("O8").Value = 0.02544529262086514
$ws.Range("Q8").Value = 0.2086513994910942
$ws.Range("R8").Value = 0.05852417302798982
$ws.Range("S8").Value = 0.3816793893129771
$ws.Range("B9").Value = 0.09523809523809523
$ws.Range("D9").Value = 0.01785714285714286
$ws.Range("F9").Value = 0.04166666666666666
$ws.Range("J9").Value = 0.1369047619047619
$ws.Range("O9").Value = 0.0119047619047619
$ws.Range("Q9").Value = 0.1964285714285714
$ws.Range("R9").Value = 0.05952380952380952
$ws.Range("S9").Value = 0.4404761904761905
$ws.Range("B10").Value = 0.11236802413273
$ws.Range("D10").Value = 0.02941176470588235
$ws.Range("E10").Value = 0.002262443438914027
$ws.Range("F10").Value = 0.07013574660633484
$ws.Range("J10").Value = 0.1372549019607843
$ws.Range("O10").Value = 0.02111613876319759
$ws.Range("Q10").Value = 0.217948717948718
$ws.Range("R10").Value = 0.05656108597285068
$ws.Range("S10").Value = 0.3529411764705883
$ws.Range("G11").Value = 0.1425
$ws.Range("J11").Value = 0.11
$ws.Range("K11").Value = 0.1925
$ws.Range("L11").Value = 0.535
$ws.Range("S11").Value = 0.02
$ws.Range("G12").Value = 0.8144796380090498
$ws.Range("J12").Value = 0.1447963800904978
$ws.Range("K12").Value = 0.01357466063348416
$ws.Range("L12").Value = 0.009049773755656109
$ws.Range("S12").Value = 0.01809954751131222
$ws.Range("G13").Value = 0.6140350877192983
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.05263157894736842
$ws.Range("F15").Value = 0.02109704641350211
$ws.Range("H15").Value = 0.1350210970464135
$ws.Range("I15").Value = 0.08438818565400844
$ws.Range("J15").Value = 0.350210970464135
$ws.Range("K15").Value = 0.0759493670886076
$ws.Range("M15").Value = 0.01265822784810127
$ws.Range("N15").Value = 0.004219409282700422
$ws.Range("O15").Value = 0.06751054852320675
$ws.Range("S15").Value = 0.2489451476793249
$ws.Range("F16").Value = 0.02450980392156863
$ws.Range("H16").Value = 0.1519607843137255
$ws.Range("I16").Value = 0.09313725490196079
$ws.Range("J16").Value = 0.3431372549019608
$ws.Range("K16").Value = 0.1813725490196078
$ws.Range("M16").Value = 0.01470588235294118
$ws.Range("O16").Value = 0.04901960784313725
$ws.Range("S16").Value = 0.142156862745098
$ws.Range("F17").Value = 0.0275049115913556
$ws.Range("H17").Value = 0.1768172888015717
$ws.Range("I17").Value = 0.06286836935166994
$ws.Range("J17").Value = 0.3713163064833006
$ws.Range("K17").Value = 0.1473477406679764
$ws.Range("M17").Value = 0.02357563850687623
$ws.Range("O17").Value = 0.07662082514734773
$ws.Range("S17").Value = 0.1139489194499018
$ws.Range("F18").Value = 0.007751937984496124
$ws.Range("H18").Value = 0.1395348837209302
$ws.Range("I18").Value = 0.08527131782945736
$ws.Range("J18").Value = 0.4496124031007752
$ws.Range("K18").Value = 0.1317829457364341
$ws.Range("M18").Value = 0.0310077519379845
$ws.Range("O18").Value = 0.04651162790697674
$ws.Range("S18").Value = 0.1085271317829457
$ws.Range("F19").Value = 0.01690821256038647
$ws.Range("H19").Value = 0.1819645732689211
$ws.Range("I19").Value = 0.07085346215780998
$ws.Range("J19").Value = 0.3735909822866345
$ws.Range("K19").Value = 0.1344605475040258
$ws.Range("M19").Value = 0.02818035426731079
$ws.Range("O19").Value = 0.07568438003220612
$ws.Range("S19").Value = 0.1167471819645733
